$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDLE TIME")

# Insert a new row before current row 15 (MARI882N.ABDELKADER) for MAKEDA.OLLIVIERRE
$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 1).Value = "MAKEDA.OLLIVIERRE"
$ws.Cells.Item(15, 2).Value = 149

# PATR5027.AMEH is now at row 20 after the insert above; update its value 162 -> 138
$ws.Cells.Item(20, 2).Value = 138

# Insert a new row before current row 22 (SEPIDEH.AZARIHASHJIN) for RARG046N.YEBOAH
$ws.Rows.Item(22).Insert()
$ws.Cells.Item(22, 1).Value = "RARG046N.YEBOAH"
$ws.Cells.Item(22, 2).Value = 159

# ZAHIDGUL.MINHAS is now at row 31 after both inserts; update its value 140 -> 120
$ws.Cells.Item(31, 2).Value = 120

# Append new row "~" / 34 at the end (row 33)
$ws.Cells.Item(33, 1).Value = "~"
$ws.Cells.Item(33, 2).Value = 34

# Update "Total Units picked by hour" sheet: make all negative values in B:E positive
$ws2 = $wb.Worksheets.Item("Total Units picked by hour")
for ($r = 2; $r -le 6; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws2.Cells.Item($r, $c)
        $v = $cell.Value()
        $cell.Value = [Math]::Abs($v)
    }
}
